$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 15-47 are entered first (in row order), and row 14 last, to reproduce
# the exact shared-string insertion order recorded by Excel in the workbook
# being reproduced (row 14's surname is a repeat of an existing string that
# was first (re)used later, while rows 15-47 introduce the brand new surname
# strings in ascending row order).
$ws.Range("A15").Value = "Bendicho Riu"
$ws.Range("B15").Value = 8
$ws.Range("C15").Value = 1924

$ws.Range("A16").Value = "Gracia Jubillà"
$ws.Range("B16").Value = 9
$ws.Range("C16").Value = 1925

$ws.Range("A17").Value = "Bendicho Jubilla"
$ws.Range("B17").Value = 10
$ws.Range("C17").Value = 1925

$ws.Range("A18").Value = "Gessé Eguanito"
$ws.Range("B18").Value = 10
$ws.Range("C18").Value = 1925

$ws.Range("A19").Value = "Ensenyat Vigo"
$ws.Range("B19").Value = 11
$ws.Range("C19").Value = 1925

$ws.Range("A20").Value = "Benavarre Gessé"
$ws.Range("B20").Value = 11
$ws.Range("C20").Value = 1925

$ws.Range("A21").Value = "Farré Ludriga"
$ws.Range("B21").Value = 12
$ws.Range("C21").Value = 1926

$ws.Range("A22").Value = "Gessé Álvarez"
$ws.Range("B22").Value = 12
$ws.Range("C22").Value = 1926

$ws.Range("A23").Value = "Benavarre Gessé"
$ws.Range("B23").Value = 13
$ws.Range("C23").Value = 1927

$ws.Range("A24").Value = "Ensenyat Vigo"
$ws.Range("B24").Value = 14
$ws.Range("C24").Value = 1927

$ws.Range("A25").Value = "Farré Ludriga"
$ws.Range("B25").Value = 14
$ws.Range("C25").Value = 1928

$ws.Range("A26").Value = "Bendicho Jubillá"
$ws.Range("B26").Value = 15
$ws.Range("C26").Value = 1929

$ws.Range("A27").Value = "Gessé Exposito"
$ws.Range("B27").Value = 15
$ws.Range("C27").Value = 1929

$ws.Range("A28").Value = "Gessé Ros"
$ws.Range("B28").Value = 16
$ws.Range("C28").Value = 1929

$ws.Range("A29").Value = "Cunyat Vigo"
$ws.Range("B29").Value = 17
$ws.Range("C29").Value = 1929

$ws.Range("A30").Value = "Ensenyat Vigo"
$ws.Range("B30").Value = 17
$ws.Range("C30").Value = 1929

$ws.Range("A31").Value = "Gessé Àlbarez"
$ws.Range("B31").Value = 18
$ws.Range("C31").Value = 1930

$ws.Range("A32").Value = "Gessé Ros"
$ws.Range("B32").Value = 19
$ws.Range("C32").Value = 1930

$ws.Range("A33").Value = "Benabarre Gessé"
$ws.Range("B33").Value = 20
$ws.Range("C33").Value = 1932

$ws.Range("A34").Value = "Gessé Ros"
$ws.Range("B34").Value = 20
$ws.Range("C34").Value = 1933

$ws.Range("A35").Value = "Fontelles Gessé"
$ws.Range("B35").Value = 21
$ws.Range("C35").Value = 1934

$ws.Range("A36").Value = "FALTA"
$ws.Range("B36").Value = 22

$ws.Range("A37").Value = "Fontelles Gessé"
$ws.Range("B37").Value = 23
$ws.Range("C37").Value = 1936

$ws.Range("A38").Value = "FALTA"
$ws.Range("B38").Value = 24

$ws.Range("A39").Value = "Puig Montanuy"
$ws.Range("B39").Value = 25
$ws.Range("C39").Value = 1938

$ws.Range("A40").Value = "Gessé Ros"
$ws.Range("B40").Value = 25
$ws.Range("C40").Value = 1939

$ws.Range("A41").Value = "Farré Tolsà"
$ws.Range("B41").Value = 26
$ws.Range("C41").Value = 1942

$ws.Range("A42").Value = "Fontelles Gessé"
$ws.Range("B42").Value = 27
$ws.Range("C42").Value = 1944

$ws.Range("A43").Value = "FALTA"
$ws.Range("B43").Value = 28

$ws.Range("A44").Value = "FALTA"
$ws.Range("B44").Value = 29

$ws.Range("A45").Value = "Fontelles Gessé"
$ws.Range("B45").Value = 30
$ws.Range("C45").Value = 1947

$ws.Range("A46").Value = "Gessé Jubillà"
$ws.Range("B46").Value = 31
$ws.Range("C46").Value = 1950

$ws.Range("A47").Value = "Fontelles Gessé"
$ws.Range("B47").Value = 32
$ws.Range("C47").Value = 1952

$ws.Range("A14").Value = "Gessé Álvarez"
$ws.Range("B14").Value = 8
$ws.Range("C14").Value = 1924

[void]$ws.Range("A48").Select()
